# ============================================================
# Generate Report for Handoff
#
# The localization CI run picked up a new handoff batch: two PNG
# screenshots (1a0e7c61-...png, 54067e3c-...png) plus the markdown
# file that used to be tracked under the old "adcad4ac-..." id is
# now tracked under a new id (8d977ceb-...). This updates the
# Overview sheet (one row per source file) and the per-locale detail
# sheets (zh-cn, de-de) with the new rows and refreshed timestamps.
# ============================================================

$wb = $excel.ActiveWorkbook

$srcBase   = "https://github.com/OpenLocalizationTest/oltest/blob/7cf06fe19a8d0c6ff8b7ed319f4b1b6ae6a9dd1a/e2e/"
$zhBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ------------------------------------------------------------
# Overview sheet
# ------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 (existing): same source, but it is now the first png and the
# handoff datetime moved on
$ws.Cells.Item(2,1).Value = "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"
$ws.Cells.Item(2,4).Value = "2016-45-18 10:45:46"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,1), ($srcBase + "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"), "", "", "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png")

# Row 3 (new): second png in the same batch
$ws.Cells.Item(3,1).Value = "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"
$ws.Cells.Item(3,2).Value = "Ready for handoff"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "2016-45-18 10:45:46"
$ws.Hyperlinks.Add($ws.Cells.Item(3,1), ($srcBase + "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"), "", "", "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png")

# Row 4 (new): the markdown file, now tracked under its new id
$ws.Cells.Item(4,1).Value = "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Cells.Item(4,2).Value = "Ready for handoff"
$ws.Cells.Item(4,3).Value = "Ready for handoff"
$ws.Cells.Item(4,4).Value = "2016-45-18 10:45:46"
$ws.Hyperlinks.Add($ws.Cells.Item(4,1), ($srcBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"), "", "", "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md")

# ------------------------------------------------------------
# zh-cn detail sheet
# ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (existing): refresh to the first png
$ws.Cells.Item(2,1).Value = "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"
$ws.Cells.Item(2,2).Value = ".png"
$ws.Cells.Item(2,4).Value = "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png"
$ws.Cells.Item(2,5).Value = "2016-03-18 10:45:41"
$ws.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(2,9).Value = "IsDependency"
$ws.Cells.Item(2,10).Value = "e2e\8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,1), ($srcBase + "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"), "", "", "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png")
$ws.Hyperlinks.Add($ws.Cells.Item(2,2), ($srcBase + "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"), "", "", ".png")
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), ($zhBase + "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png"), "", "", "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png")

# Row 3 (new): second png
$ws.Cells.Item(3,1).Value = "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"
$ws.Cells.Item(3,2).Value = ".png"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "92f96ea6697ad4094c3d741c5537e4890c43102c.png"
$ws.Cells.Item(3,5).Value = "2016-03-18 10:45:41"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "IsDependency"
$ws.Cells.Item(3,10).Value = "e2e\8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Hyperlinks.Add($ws.Cells.Item(3,1), ($srcBase + "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"), "", "", "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png")
$ws.Hyperlinks.Add($ws.Cells.Item(3,2), ($srcBase + "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"), "", "", ".png")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), ($zhBase + "92f96ea6697ad4094c3d741c5537e4890c43102c.png"), "", "", "92f96ea6697ad4094c3d741c5537e4890c43102c.png")

# Row 4 (new): the markdown file
$ws.Cells.Item(4,1).Value = "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Cells.Item(4,2).Value = ".md"
$ws.Cells.Item(4,3).Value = "Ready for handoff"
$ws.Cells.Item(4,4).Value = "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.zh-cn.xlf"
$ws.Cells.Item(4,5).Value = "2016-03-18 10:45:41"
$ws.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,9).Value = "Include"
$ws.Hyperlinks.Add($ws.Cells.Item(4,1), ($srcBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"), "", "", "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,2), ($srcBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"), "", "", ".md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4), ($zhBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.zh-cn.xlf"), "", "", "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.zh-cn.xlf")

# ------------------------------------------------------------
# de-de detail sheet
# ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (existing): refresh to the first png
$ws.Cells.Item(2,1).Value = "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"
$ws.Cells.Item(2,2).Value = ".png"
$ws.Cells.Item(2,4).Value = "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png"
$ws.Cells.Item(2,5).Value = "2016-03-18 10:45:46"
$ws.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(2,9).Value = "IsDependency"
$ws.Cells.Item(2,10).Value = "e2e\8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,1), ($srcBase + "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"), "", "", "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png")
$ws.Hyperlinks.Add($ws.Cells.Item(2,2), ($srcBase + "1a0e7c61-e871-4d64-83d0-1a75907fba1f.png"), "", "", ".png")
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), ($deBase + "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png"), "", "", "4352a99d9ba946de2ddfc0bcc3aee89b2c209225.png")

# Row 3 (new): second png
$ws.Cells.Item(3,1).Value = "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"
$ws.Cells.Item(3,2).Value = ".png"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "92f96ea6697ad4094c3d741c5537e4890c43102c.png"
$ws.Cells.Item(3,5).Value = "2016-03-18 10:45:46"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "IsDependency"
$ws.Cells.Item(3,10).Value = "e2e\8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Hyperlinks.Add($ws.Cells.Item(3,1), ($srcBase + "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"), "", "", "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png")
$ws.Hyperlinks.Add($ws.Cells.Item(3,2), ($srcBase + "54067e3c-3cb7-4d31-9cf2-c0832c79ca99.png"), "", "", ".png")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), ($deBase + "92f96ea6697ad4094c3d741c5537e4890c43102c.png"), "", "", "92f96ea6697ad4094c3d741c5537e4890c43102c.png")

# Row 4 (new): the markdown file
$ws.Cells.Item(4,1).Value = "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"
$ws.Cells.Item(4,2).Value = ".md"
$ws.Cells.Item(4,3).Value = "Ready for handoff"
$ws.Cells.Item(4,4).Value = "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.de-de.xlf"
$ws.Cells.Item(4,5).Value = "2016-03-18 10:45:46"
$ws.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,9).Value = "Include"
$ws.Hyperlinks.Add($ws.Cells.Item(4,1), ($srcBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"), "", "", "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,2), ($srcBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.md"), "", "", ".md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4), ($deBase + "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.de-de.xlf"), "", "", "8d977ceb-2f5b-4137-b1e1-372ff0081a1d.6b0511bd435c0d733d9396980768972d647fbf7e.de-de.xlf")

Write-Host "Handoff report generated: Overview/zh-cn/de-de rows 3-4 added, row 2 refreshed."
